# Update Negative_Manifest sheet: append 25 new rows (r=40..64) of negative,
# no_meltpatch training-subject records, and move the active selection /
# view down to the freshly appended block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('n39', 'n39_IMG_3174.jpeg', 'True', 'no_meltpatch', 'negative'),
    @('n40', 'n40_IMG_3174HorFlip.jpeg', 'True', 'no_meltpatch', 'negative'),
    @('n41', 'n41_IMG_3174HorVertFlip.jpeg', 'True', 'no_meltpatch', 'negative'),
    @('n42', 'n42_IMG_3174VertFlip.jpeg', 'True', 'no_meltpatch', 'negative'),
    @('n43', 'n43_IMG_3176.jpeg', 'True', 'no_meltpatch', 'negative'),
    @('n44', 'n44_IMG_3176HorFlip.jpeg', 'True', 'no_meltpatch', 'negative'),
    @('n45', 'n45_IMG_3176HorVertFlip.jpeg', 'True', 'no_meltpatch', 'negative'),
    @('n46', 'n46_IMG_3176VertFlip.jpeg', 'True', 'no_meltpatch', 'negative'),
    @('n47', 'n47_IMG_3177.jpeg', 'True', 'no_meltpatch', 'negative'),
    @('n48', 'n48_IMG_3177HorFlip.jpeg', 'True', 'no_meltpatch', 'negative'),
    @('n49', 'n49_IMG_3177HorVertFlip.jpeg', 'True', 'no_meltpatch', 'negative'),
    @('n50', 'n50_IMG_3177VertFlip.jpeg', 'True', 'no_meltpatch', 'negative'),
    @('n51', 'n51_IMG_3178.jpeg', 'True', 'no_meltpatch', 'negative'),
    @('n52', 'n52_IMG_3070.jpeg', 'True', 'no_meltpatch', 'negative'),
    @('n53', 'n53_IMG_3070HorFlip.jpeg', 'True', 'no_meltpatch', 'negative'),
    @('n54', 'n54_IMG_3070HorVertFlip.jpeg', 'True', 'no_meltpatch', 'negative'),
    @('n55', 'n55_IMG_3070VertFlip.jpeg', 'True', 'no_meltpatch', 'negative'),
    @('n56', 'n56_IMG_3072.jpeg', 'True', 'no_meltpatch', 'negative'),
    @('n57', 'n57_IMG_3072HorFlip.jpeg', 'True', 'no_meltpatch', 'negative'),
    @('n58', 'n58_IMG_3072HorVertFlip.jpeg', 'True', 'no_meltpatch', 'negative'),
    @('n59', 'n59_IMG_3072VertFlip.jpeg', 'True', 'no_meltpatch', 'negative'),
    @('n60', 'n60_IMG_3073.jpeg', 'True', 'no_meltpatch', 'negative'),
    @('n61', 'n61_IMG_3073HorFlip.jpeg', 'True', 'no_meltpatch', 'negative'),
    @('n62', 'n62_IMG_3073HorVertFlip.jpeg', 'True', 'no_meltpatch', 'negative'),
    @('n63', 'n63_IMG_3073VertFlip.jpeg', 'True', 'no_meltpatch', 'negative'),
)

$startRow = 40
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    # Column C holds the literal text "True" (not a boolean) in this sheet,
    # so force text entry the way Excel's UI does (leading apostrophe),
    # then strip the resulting quote-prefix cell style so the cell matches
    # the plain, unstyled text cells used throughout the rest of the sheet.
    $ws.Cells.Item($row, 3).Value = "'" + $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
}

$endRow = $startRow + $data.Count - 1
$ws.Range("C$startRow`:C$endRow").ClearFormats()

# Move the view/selection onto the newly appended block, like the source
# workbook did after the append (top row scrolled to A35, selection on the
# second half of the freshly-added rows).
[void]$ws.Range("A40:E52").Select()
